$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.044366416687426, 1.051754987288968, 1.051835858366668, 1.061815470365856, 1.02359499962809, 1.049431886253281, 1.05450561112702, 1.054586258161316, 1.064538506498105, 1.050922199903281),
    @(1.02, 1.045742149378772, 1.05300963872044, 1.053073294039902, 1.063166637302955, 1.023504579208684, 1.050452735914316, 1.055571836503854, 1.055635328350022, 1.065703017744349, 1.051944499287927),
    @(1.02, 1.046632404330529, 1.053821801990218, 1.05387431091859, 1.064041456525528, 1.023443335729026, 1.051112912620486, 1.056261517196469, 1.05631389812942, 1.066456516076861, 1.052605613520744),
    @(1.02, 1.047006686586925, 1.054163314693474, 1.054211135995032, 1.06440935887725, 1.023416932628374, 1.051390362047623, 1.056551404318079, 1.056599111670253, 1.066773285011319, 1.052883456957891),
    @(1.02, 1.04706953145298, 1.054220660905101, 1.054267695032026, 1.064471138912839, 1.023412460900047, 1.051436941886826, 1.056600074507989, 1.056646996993568, 1.066826471881925, 1.052930102945808),
    @(1.02, 1.046637405428699, 1.05382636498488, 1.053878811283686, 1.064046371946768, 1.023442985510381, 1.051116620259462, 1.05626539089405, 1.056317709388508, 1.066460748762541, 1.052609326424992),
    @(1.02, 1.044831341031003, 1.052178937991249, 1.052253992242456, 1.062271996856341, 1.023565007130014, 1.049776967580925, 1.054865997819324, 1.054940849278644, 1.064932063894522, 1.051267771285991),
    @(1.02, 1.041649145447721, 1.049278291101655, 1.049393136725764, 1.059149201931629, 1.023759139049229, 1.047413299785727, 1.052398153224373, 1.052512634553902, 1.062238068275243, 1.048900746811407),
    @(1.02, 1.039527660251703, 1.047345923567211, 1.047487271135893, 1.057069762400187, 1.023874603872166, 1.045835332076313, 1.050751468326264, 1.050892322931897, 1.060441732524758, 1.04732053820685),
    @(1.02, 1.038608975565838, 1.046509480238737, 1.046662298636129, 1.056169875254395, 1.023921309042974, 1.04515150093701, 1.050038055707305, 1.050190320172058, 1.0596637816879, 1.046635735948893),
    @(1.02, 1.038267720380208, 1.046198827334765, 1.046355905913463, 1.055835692144694, 1.023938164268129, 1.04489740832854, 1.049773001264149, 1.04992950245072, 1.059374793956056, 1.046381282499774),
    @(1.02, 1.03834092151968, 1.046265461708765, 1.046421626489291, 1.055907372219209, 1.023934571058996, 1.044951916053831, 1.049829859162049, 1.049985451619735, 1.059436783860783, 1.046435867632287),
    @(1.02, 1.038580767634265, 1.046483800788705, 1.046636971345522, 1.056142250064666, 1.023922712353276, 1.045130499360716, 1.050016147491192, 1.050168762192915, 1.059639894342762, 1.046614704547952),
    @(1.02, 1.038728542642933, 1.046618331820899, 1.046769657462063, 1.056286975891407, 1.023915340507297, 1.045240518774254, 1.050130917673571, 1.050281697497583, 1.059765034380409, 1.046724880201676),
    @(1.02, 1.039588628509276, 1.047401441264887, 1.047542027462983, 1.057129495671544, 1.023871434968932, 1.04588070362447, 1.050798806753474, 1.050938903922176, 1.060493359633399, 1.047365974187811),
    @(1.02, 1.04012811628432, 1.047892738944827, 1.048026587395859, 1.057658123746516, 1.023843013862359, 1.04628212240886, 1.051217650106968, 1.05135104372676, 1.060950182879899, 1.047767963032918),
    @(1.02, 1.040442784168372, 1.048179332254409, 1.048309250231855, 1.057966514356165, 1.023826118601784, 1.046516209300371, 1.051461917198398, 1.051591399786307, 1.06121662825493, 1.048002382354664),
    @(1.02, 1.040550076871775, 1.048277057963751, 1.048405635670002, 1.058071676310929, 1.023820303822167, 1.046596017831052, 1.051545199754905, 1.051673348615697, 1.061307477323299, 1.048082304222613),
    @(1.02, 1.040070235028804, 1.047840024515355, 1.04797459590962, 1.05760140173062, 1.023846096030192, 1.046239059569861, 1.051172716044948, 1.051306828991292, 1.060901171324278, 1.047724839039747),
    @(1.02, 1.038510139377035, 1.046419504311286, 1.046573556607911, 1.056073082351267, 1.02392621804881, 1.045077913461566, 1.049961291936291, 1.050114783551098, 1.059580084006143, 1.046562043970793),
    @(1.02, 1.037529157455005, 1.045526593110581, 1.045692890434128, 1.055112599656355, 1.023973741382756, 1.044347348358118, 1.049199265728368, 1.0493649333783, 1.058749334174753, 1.045830441381102),
    @(1.02, 1.038049204410176, 1.045999921912478, 1.046159728200397, 1.055621729948876, 1.023948818288664, 1.044734683763612, 1.049603264987206, 1.049762478629855, 1.059189743779694, 1.046218326847299),
    @(1.02, 1.040096389095662, 1.047863843811324, 1.04799808854009, 1.057627031810961, 1.023844704313072, 1.046258517979599, 1.051193019928095, 1.051326807844722, 1.060923317573232, 1.04774432508266),
    @(1.02, 1.042471805672414, 1.050027919167128, 1.050132482710052, 1.059956078459364, 1.02371141742819, 1.048024739200764, 1.053036395380807, 1.053140639990816, 1.062934578261755, 1.049513054540547)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $cols[$c]).Value = $rowVals[$c]
    }
}
